$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header columns for team record: Wins, Losses, Ties
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy header style (bold, centered, bordered) from an existing header cell
$ws.Range("AC1").Copy() | Out-Null
$ws.Range("AD1:AF1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

# Re-apply header values, since paste special formats only shouldn't touch them,
# but ensure values are still correct
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

$lastRow = 48
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 75   # AD
    $ws.Cells.Item($r, 31).Value = 86   # AE
    $ws.Cells.Item($r, 32).Value = 0    # AF
}
